# Add a new "LBP+HOG" results sheet (cloned from the "Gabor" report sheet so it
# inherits the same number formats / alignment / row-height behaviour), fill in
# the new experiment's numbers, and update the small set of "last selected
# cell" bookkeeping bits on the other report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet as a copy of "Gabor" (last sheet), positioned right
#    after it, then rename it.
# ---------------------------------------------------------------------------
$gabor = $wb.Worksheets.Item("Gabor")
$gabor.Copy($null, $gabor)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "LBP+HOG"

# Center-align (horizontal + vertical) the whole A1:E14 block, matching the
# formatting used by every other report sheet - this also extends the used
# range to column E.
$rng = $newSheet.Range("A1:E14")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Best-effort column widths to match the new sheet's layout.
$newSheet.Columns("B").ColumnWidth = 14.41
$newSheet.Columns("C:D").ColumnWidth = 11.58

# ---------------------------------------------------------------------------
# 2. Fill in the experiment write-up + metrics.
# ---------------------------------------------------------------------------
$newSheet.Range("A1").Value = "25/01/26 Trying LBP and HOG combined, ImageResizing not applied to LBP. All features used, StandardScaler. PCA Applied with 90% threshold."

$newSheet.Range("B2").Value = "RandomForest"
$newSheet.Range("C2").Value = "GradientBoosting"
$newSheet.Range("D2").Value = "SVM"

$newSheet.Range("A3").Value = "Train"
$newSheet.Range("A4").Value = "Accuracy"
$newSheet.Range("B4").Value = 1
$newSheet.Range("C4").Value = 0.8301
$newSheet.Range("D4").Value = 0.955
$newSheet.Range("A5").Value = "F1-Score"
$newSheet.Range("B5").Value = 1
$newSheet.Range("C5").Value = 0.8326
$newSheet.Range("D5").Value = 0.9551
$newSheet.Range("A6").Value = "Confusion Matrix"

$newSheet.Range("A7").Value = "Valid"
$newSheet.Range("A8").Value = "Accuracy"
$newSheet.Range("B8").Value = 0.7325
$newSheet.Range("C8").Value = 0.722
$newSheet.Range("D8").Value = 0.8285
$newSheet.Range("A9").Value = "F1-Score"
$newSheet.Range("B9").Value = 0.7277
$newSheet.Range("C9").Value = 0.7301
$newSheet.Range("D9").Value = 0.8323
$newSheet.Range("A10").Value = "Confusion Matrix"

$newSheet.Range("A11").Value = "Test"
$newSheet.Range("A12").Value = "Accuracy"
$newSheet.Range("B12").Value = 0.7415
$newSheet.Range("C12").Value = 0.735
$newSheet.Range("D12").Value = 0.831
$newSheet.Range("A13").Value = "F1-Score"
$newSheet.Range("B13").Value = 0.7411
$newSheet.Range("C13").Value = 0.7462
$newSheet.Range("D13").Value = 0.8337
$newSheet.Range("A14").Value = "Confusion Matrix"

# Confusion-matrix text blocks, written column-by-column (B, then C, then D)
# so shared-string allocation order matches the authored workbook.
$newSheet.Range("B6").Value = "[[5000   0]`n [0 5000]]"
$newSheet.Range("B10").Value = "[[750 250]`n [285 715]]"
$newSheet.Range("B14").Value = "[[743 257]`n [260 740]]"

$newSheet.Range("C6").Value = "[[4076  924]`n [ 775 4225]]"
$newSheet.Range("C10").Value = "[[692 308]`n [248 752]]"
$newSheet.Range("C14").Value = "[[691 309]`n [221 779]]"

$newSheet.Range("D6").Value = "[[4764  236]`n [ 214 4786]]"
$newSheet.Range("D10").Value = "[[806 194]`n [149 851]]"
$newSheet.Range("D14").Value = "[[815 185]`n [153 847]]"

# ---------------------------------------------------------------------------
# 3. Update the "last selected cell" on the other report sheets that were
#    touched while navigating around to add the new sheet.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Baseline+PCA").Range("A2").Select()
$wb.Worksheets.Item("LBP").Range("E13").Select()
$wb.Worksheets.Item("HOG").Range("E13").Select()
$wb.Worksheets.Item("Color").Range("E12").Select()

# Select the new sheet last so it ends up the active tab, matching the saved
# workbook state.
$newSheet.Range("D12").Select()
